$wb = $excel.ActiveWorkbook

# --- "Metadata" sheet (Property / Value table) ---
$meta = $wb.Worksheets.Item("Metadata")

# Version 5.0.0 -> 6.0.0
$meta.Range("B3").Value = "6.0.0"

# Date updated
$meta.Range("B8").Value = "2022-01-21T20:46:54+00:00"

# Publisher value was blank -> "Alvearie Team"
$meta.Range("B9").Value = "Alvearie Team"

# Row 10 was "Contact" / "No display for ContactDetail" -> becomes "Jurisdiction" / "United States of America"
$meta.Range("A10").Value = "Jurisdiction"
$meta.Range("B10").Value = "United States of America"

# Row 11 was a duplicate "Contact" / "No display for ContactDetail" row -> remove it entirely,
# shifting the remaining rows (Description, Purpose, Copyright, ...) up by one.
$meta.Rows.Item(11).Delete()

# --- "Elements" sheet (big FHIR element grid) ---
$elem = $wb.Worksheets.Item("Elements")

# Root "Extension" row: Short/Definition now reflect the actual extension's purpose.
$elem.Range("K2").Value = "Confidence Score"
$elem.Range("L2").Value = "Confidence score for the match"
